$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Oklo Inc. / OKLO)
$ws.Range("D2").Value = 109.99
$ws.Range("E2").Value = 56.3
$ws.Range("F2").Value = 23.97
$ws.Range("I2").Value = 70
$ws.Range("K2").Value = 59.1
$ws.Range("N2").Value = 53.62998959737769

# Row 3 (NuScale Power Corporation / SMR)
$ws.Range("D3").Value = 22.68
$ws.Range("E3").Value = 48.5
$ws.Range("F3").Value = 19.12
$ws.Range("K3").Value = 55.1
$ws.Range("N3").Value = 53.62998959737769
